# ===================================================================
# Section A: Paragraph 44 "Mantenimiento" intro + two new List Bullet items
# ===================================================================
$d = $word.ActiveDocument

$p44 = $d.Paragraphs.Item(44)
$p44.Range.Text = 'Para añadir nuevos campos de extracción o ajustar el filtrado:'

$p44.Range.InsertParagraphAfter()
$p45 = $d.Paragraphs.Item(45)
$p45.Range.Text = 'Patrones: Añade el patrón regex en el diccionario `self.patterns` en `extractor.py`.'
$p45.Style = "List Bullet"

$p45.Range.InsertParagraphAfter()
$p46 = $d.Paragraphs.Item(46)
$p46.Range.Text = 'Reglas de exclusión: Modifica el método `extraer_datos_pdf` en `extractor.py`.'
$p46.Style = "List Bullet"

# ===================================================================
# Section B: Paragraph 41 "GitHub" bullet -> ".gitignore" bullet
# ===================================================================
$p41 = $d.Paragraphs.Item(41)
$p41.Range.Text = 'Los informes generados y los datos de entrada están excluidos en el `.gitignore`.'

# ===================================================================
# Section C: New "Filtrado Inteligente" block inserted before paragraph 37
#            (the blank paragraph that precedes "Despliegue y GitHub")
# ===================================================================
$anchor = $d.Paragraphs.Item(37)
for ($k = 0; $k -lt 5; $k++) {
    $anchor.Range.InsertParagraphBefore()
}
# Paragraph 37 is now a brand new blank separator paragraph (left blank on purpose)

$h = $d.Paragraphs.Item(38)
$h.Range.Text = '🛡️ Lógica de Filtrado Inteligente'
$h.Style = "Heading 1"

$intro = $d.Paragraphs.Item(39)
$intro.Range.Text = 'Para asegurar que los informes contengan solo facturas válidas, se aplican dos niveles de filtrado:'

$f1 = $d.Paragraphs.Item(40)
$f1.Range.Text = '**Filtro por Nombre**: Se omiten archivos que contengan palabras como "CONTRATO" o "CARTA" en su nombre.'
$f1.Style = "List Number"

$f2 = $d.Paragraphs.Item(41)
$f2.Range.Text = '**Filtro por Contenido**: Si tras procesar el PDF no se encuentra un Número de Factura Y el Total es 0, el documento se considera irrelevante y no se añade al Excel.'
$f2.Style = "List Number"

# ===================================================================
# Section D: Paragraph 35 "Número de Factura" bullet + two new bullets
#            ("Fechas" and "Importes")
# ===================================================================
$p35 = $d.Paragraphs.Item(35)
$p35.Range.Text = '**Número de Factura**: Busca palabras clave como "Factura nº", "Nº Factura", "Invoice". Se ha refinado para ignorar etiquetas como "Tlfno" o "Fax".'

$p35.Range.InsertParagraphAfter()
$p36 = $d.Paragraphs.Item(36)
$p36.Range.Text = '**Fechas**: Se extraen por separado la "Fecha de Factura" (emisión) y la "Fecha de Cargo" (vencimiento/cobro).'

$p36.Range.InsertParagraphAfter()
$p37b = $d.Paragraphs.Item(37)
$p37b.Range.Text = '**Importes**: Captura tanto la "Base Imponible" como el "Total" de la factura.'
